$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header band: switch to left alignment ---
# C4 alone first so its existing style slot (title font+fill) is updated in place.
$ws.Range("C4").HorizontalAlignment = -4131
# Rest of row 4 and row 5 (dark-fill band) also become left aligned (new shared style).
$ws.Range("D4:I5").HorizontalAlignment = -4131
$ws.Range("C5").HorizontalAlignment = -4131

# --- New task rows (7-13) : text/number values ---
$ws.Range("C7").Value = "Entregar novo protótipo de Calculadora"
$ws.Range("D7").Value = "Alta"
$ws.Range("E7").Value = "Fazendo"
$ws.Range("G7").Value = "Enrico e Samuel"

$ws.Range("C8").Value = "Modificações na documentação"
$ws.Range("D8").Value = "Alta"
$ws.Range("E8").Value = "Fazendo"
$ws.Range("G8").Value = "Henry"

$ws.Range("C9").Value = "Novo Banco de Dados com DER"
$ws.Range("D9").Value = "Alta"
$ws.Range("E9").Value = "Fazendo"
$ws.Range("G9").Value = "Cristhian e Kaue"

$ws.Range("C10").Value = "Modificações nos Slides de Apresentação"
$ws.Range("D10").Value = "Media"
$ws.Range("E10").Value = "Fazendo"
$ws.Range("G10").Value = "Larissa e Tabata"

$ws.Range("C11").Value = "Protótipo do Site no Figma"
$ws.Range("D11").Value = "Media"
$ws.Range("E11").Value = "Fazendo"
$ws.Range("G11").Value = "Kaue"

$ws.Range("C12").Value = "Verificação do código de arduino"
$ws.Range("D12").Value = "Alta"
$ws.Range("E12").Value = "Fazendo"
$ws.Range("G12").Value = "Cristhian e Kaue"

$ws.Range("C13").Value = "Criação do Site por código"
$ws.Range("D13").Value = "Media"
$ws.Range("E13").Value = "A Fazer"

# --- Prazo (date) column: set values then apply the built-in short-date format ---
$ws.Range("F7").Value = 45394
$ws.Range("F8").Value = 45394
$ws.Range("F9").Value = 45394
$ws.Range("F10").Value = 45394
$ws.Range("F11").Value = 45396
$ws.Range("F12").Value = 45394
$ws.Range("F13").Value = 45398

# Format F7 (light-blue row) and clone that format onto the other light rows.
$ws.Range("F7").NumberFormat = "mm-dd-yy"
$ws.Range("F7").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F13").PasteSpecial(-4122)

# Format F8 (gray row) and clone that format onto the other gray rows.
$ws.Range("F8").NumberFormat = "mm-dd-yy"
$ws.Range("F8").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Range("F12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Column C widened to fit the longer task names ---
$ws.Columns.Item(3).ColumnWidth = 36

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection cursor moved ---
$ws.Range("I13").Select()
